$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated/regrouped antibiotic classes and extended discovery timeline.
# (class, year as Excel date serial, type)
$data = @(
    @("Arsphenamines",               2558),
    @("β-Lactams",                   10594),
    @("Sulphonamides",               11689),
    @("Polypeptides",                14246),
    @("Aminoglycosides",             16072),
    @("Nitrofurans",                 16438),
    @("Tetracyclines",               17533),
    @("Phenicols",                   17168),
    @("Macrolides",                  18994),
    @("Streptogramins",              19360),
    @("Glycopeptides",               19725),
    @("Azoles",                      21551),
    @("Quinolones and Lincosamides", 22647),
    @("Oxazolidinones",              31778),
    @("Diarylquinolines",            37987)
)

$ws.Range("A1").Value = "class"
$ws.Range("B1").Value = "year"
$ws.Range("C1").Value = "type"

$row = 2
foreach ($entry in $data) {
    if ($row -gt 14) {
        # New rows: copy the date-formatted style from the row above
        # before writing, so the new cells pick up the same number format.
        $ws.Cells.Item($row - 1, 2).Copy() | Out-Null
        $ws.Cells.Item($row, 2).PasteSpecial(-4122) | Out-Null
    }
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value2 = $entry[1]
    $ws.Cells.Item($row, 3).Value = "timeline"
    $row++
}

$excel.CutCopyMode = 0

$ws.Application.Goto($ws.Range("A3"), $true)
$ws.Range("D7").Select()
